# Daily attendance processing - 2025-11-28 05:52:28
# Normalize the "Recorded By" (column G) cell values: for any cell that lists
# more than one recorder separated by ", ", swap the first two entries while
# leaving any further entries (and single-entry cells) untouched. Cells whose
# value is exactly "System, backup@backdoor.com" are left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $current = $cell.Text

    if ($current -eq $null -or $current -eq "") {
        continue
    }

    if ($current -eq "System, backup@backdoor.com") {
        continue
    }

    $parts = $current -split ", "
    if ($parts.Count -ge 2) {
        $swapped = @($parts[1], $parts[0])
        if ($parts.Count -gt 2) {
            $swapped += $parts[2..($parts.Count - 1)]
        }
        $newValue = [string]::Join(", ", $swapped)
        $cell.Value = $newValue
    }
}
